$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Controleur CA": add row 27 (new "Couche prises" control), matching
# the pattern already used by row 26 ("Couche point contrôle"). Also the
# header row no longer needs its extra-tall custom height.
# ---------------------------------------------------------------------------
$wsCtrl = $wb.Worksheets.Item("Controleur CA")

# Duplicate row 26's formatting into row 27, then overwrite with new values.
$wsCtrl.Range("A26:N26").Copy($wsCtrl.Range("A27:N27"))

$wsCtrl.Range("A27").Value = 27
$wsCtrl.Range("B27").Value = "commande d'accès"
$wsCtrl.Range("C27").Value = "Structuration des couches"
$wsCtrl.Range("D27").Value = "Couche prises"
$wsCtrl.Range("E27").Value = "[vide]"
$wsCtrl.Range("F27").Value = "Vérifier que la liste des champs de la couche prises correspond aux spécifications QGIS"
$wsCtrl.Range("G27").Value = "[vide]"
$wsCtrl.Range("H27").Value = "La structuration des champs de la couche prises est incorrecte"
$wsCtrl.Range("I27").Value = "Majeure"

# ---------------------------------------------------------------------------
# Sheet "CD21": row 26 gains a taller row height, and row 27 (previously an
# almost-empty placeholder row) is filled in with formulas mirroring the
# "Controleur CA" row 27 content.
# ---------------------------------------------------------------------------
$wsCD21 = $wb.Worksheets.Item("CD21")
$wsCD21.Rows.Item(26).RowHeight = 30

$wsCD21.Range("A26:M26").Copy($wsCD21.Range("A27:M27"))
$wsCD21.Range("A27").Value = 27
$wsCD21.Range("B27").Formula = "='Controleur CA'!B27"
$wsCD21.Range("C27").Formula = "='Controleur CA'!F27"
$wsCD21.Rows.Item(27).RowHeight = 30

$wsCD21.Activate()
$wsCD21.Range("C29").Select()

# ---------------------------------------------------------------------------
# Sheets "CD39", "CD58", "CD70", "CD71", "testv1", "testv2": append a brand
# new row 27, mirroring row 26's formatting/formulas.
# ---------------------------------------------------------------------------
$otherSheets = @("CD39", "CD58", "CD70", "CD71", "testv1", "testv2")
foreach ($name in $otherSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A26:M26").Copy($ws.Range("A27:M27"))
    $ws.Range("A27").Value = 27
    $ws.Range("B27").Formula = "='Controleur CA'!B27"
    $ws.Range("C27").Formula = "='Controleur CA'!F27"
    $ws.Rows.Item(27).RowHeight = 30

    $ws.Activate()
    $ws.Rows.Item(27).Select()
}

# testv1's window used to be scrolled so column D showed first; restore to A.
$wsTestv1 = $wb.Worksheets.Item("testv1")
$wsTestv1.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$wsTestv1.Rows.Item(27).Select()

# ---------------------------------------------------------------------------
# Re-activate "Controleur CA" last so it remains the selected tab, and park
# the selection on H28 (just below the newly-added row).
# ---------------------------------------------------------------------------
$wsCtrl.Activate()
$wsCtrl.Range("H28").Select()
